# Update to CDM online manual Appendix B
# Rewrites the "Raw"/"AED" data table on Sheet3 with revised mean ± SE
# text values (no longer computed via formulas), restyles the Light/Dark
# data columns (E:F) as Text-formatted cells, and moves the active
# selection to H2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# --- Replace formula-derived / placeholder values with literal text ---
# Column E = "Light", Column F = "Dark" (headers in row 1 are unchanged)

# AMM
$ws.Range("E2").Value = "253 ± 88.4"
$ws.Range("F2").Value = "172.5 ± 107.6"
$ws.Range("E3").Value = "6.07 ± 2.12"
$ws.Range("F3").Value = "4.14 ± 2.58"

# NIT
$ws.Range("E4").Value = "267.2 ± 269.8"
$ws.Range("F4").Value = "2.90 ± 43.6"
$ws.Range("E5").Value = "6.41 ± 6.47"
$ws.Range("F5").Value = "0.07 ± 1.05"

# FRP
$ws.Range("E6").Value = "5.6 ± 5.8"
$ws.Range("F6").Value = "28.3 ± 16.9"
$ws.Range("E7").Value = "0.13 ± 0.14"
$ws.Range("F7").Value = "0.68 ± 0.41"

# OXY
$ws.Range("E8").Value = "1546 ± 2857"
$ws.Range("F8").Value = "-4925 ± 541 "
$ws.Range("E9").Value = "37.1 ± 68.57"
$ws.Range("F9").Value = "-118.2  ± 12.98"

# --- Re-style the Light/Dark data range as Text-formatted cells ---
$ws.Range("E2:E9").NumberFormat = "@"
$ws.Range("F2:F9").NumberFormat = "@"

# --- Move the active selection on Sheet3 to H2 ---
$ws.Range("H2").Select()
